## Edit commit: "Thu, May 07, 2020  1:07:04 AM"
##
## This commit:
##  1. Changes the table style (on the "Sources of finance" table) to the
##     built-in style {2192E542-0AC5-4C6D-BE96-00616F948946}.
##  2. Swaps the two themes in the package: the "Integral" design theme
##     (theme1.xml, used by the slide master / the presentation design)
##     and the "Office Theme" (theme2.xml, used only by the notes master)
##     trade places, i.e. the slide master's theme colours become the
##     stock "Office" palette, and (in the source OOXML) the notes
##     master's theme becomes "Integral".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table with the new built-in table style id.
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{2192E542-0AC5-4C6D-BE96-00616F948946}")
        }
    }
}

# ---------------------------------------------------------------------
# 2) Re-colour the presentation's theme (theme1.xml) so that it carries
#    the standard Office theme palette (what the commit moves onto the
#    slide master side of the swap).
# ---------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme colours, VBA RGB() long values (0xBBGGRR) for:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$themeColors.Colors(1).RGB  = 0            # dk1      000000
$themeColors.Colors(2).RGB  = 16777215     # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388      # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391     # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939     # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501      # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845     # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407        # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308     # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456      # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797     # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477      # folHlink 954F72
